$wb = $excel.ActiveWorkbook

# --- Overview sheet: handoff status text changed ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- zh-cn sheet: handoff generation failed, so clear the handoff file link
#     and reset the handoff datetime; flip reason from Include to Ignored ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$2') {
        $h.Delete()
    }
}
$wsZhCn.Range("C2").Clear()
$wsZhCn.Range("D2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Ignored"

# --- de-de sheet: same treatment ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$2') {
        $h.Delete()
    }
}
$wsDeDe.Range("C2").Clear()
$wsDeDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Ignored"
